$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 80
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 76
$ws.Range("F2").Value = 4

# Update row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 20

# Update row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 22
$ws.Range("F4").Value = 8

# Remove rows 5 and 6 (no longer part of the data range)
$ws.Range("A5:F6").Delete()
